# Scheduled-runner update: refresh crafting-leve profit figures (currentAveragePrice /
# currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ columns H..N) per job sheet,
# reflecting newly pulled market-board prices. A few rows whose items have no current
# market data had their price/profit cells cleared instead of recalculated.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (ALC)
$ws.Range("H2").Value = 227.5
$ws.Range("I2").Value = 234.28572
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 234.28572
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = -121.28572
$ws.Range("N2").Value = -406

# Row 87 (ALC)
$ws.Range("H87").Value = 17230.477
$ws.Range("J87").Value = 17230.477
$ws.Range("L87").Value = 17230.477
$ws.Range("N87").Value = -19726.477

# Row 90 (ALC)
$ws.Range("H90").Value = 17230.477
$ws.Range("J90").Value = 17230.477
$ws.Range("L90").Value = 51691.431
$ws.Range("N90").Value = -64171.431

# Row 107 (ALC)
$ws.Range("H107").Value = 426.2
$ws.Range("I107").Value = 226.78947
$ws.Range("K107").Value = 226.78947
$ws.Range("M107").Value = 1693.21053

# Row 129 (ALC)
$ws.Range("H129").Value = 811.87933
$ws.Range("J129").Value = 883.1667
$ws.Range("L129").Value = 2649.5001
$ws.Range("N129").Value = -12649.5001

# Row 138 (ALC)
$ws.Range("H138").Value = 2176.5
$ws.Range("I138").Value = 1663
$ws.Range("J138").Value = 2545.8596
$ws.Range("K138").Value = 4989
$ws.Range("L138").Value = 7637.578799999999
$ws.Range("M138").Value = 151
$ws.Range("N138").Value = -17917.5788

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 7230.722
$ws.Range("I32").Value = 5892.256
$ws.Range("J32").Value = 20950
$ws.Range("K32").Value = 5892.256
$ws.Range("L32").Value = 20950
$ws.Range("M32").Value = -5605.256
$ws.Range("N32").Value = -21524

# Row 34 (ARM)
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

# Row 37 (ARM)
$ws.Range("H37").Value = 11818.182
$ws.Range("I37").Value = 11818.182
$ws.Range("K37").Value = 11818.182
$ws.Range("M37").Value = -11545.182

# Row 61 (ARM)
$ws.Range("H61").Value = 2247.0862
$ws.Range("I61").Value = 2107.9375
$ws.Range("J61").Value = 2915
$ws.Range("K61").Value = 2107.9375
$ws.Range("L61").Value = 2915
$ws.Range("M61").Value = -1895.9375
$ws.Range("N61").Value = -3339

# Row 136 (ARM)
$ws.Range("H136").Value = 2247.0862
$ws.Range("I136").Value = 2107.9375
$ws.Range("J136").Value = 2915
$ws.Range("K136").Value = 6323.8125
$ws.Range("L136").Value = 8745
$ws.Range("M136").Value = -3773.8125
$ws.Range("N136").Value = -13845

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (BSM)
$ws.Range("H22").Value = 233.2
$ws.Range("I22").Value = 104.57143
$ws.Range("J22").Value = 533.3333
$ws.Range("K22").Value = 104.57143
$ws.Range("L22").Value = 533.3333
$ws.Range("M22").Value = 68.42856999999999
$ws.Range("N22").Value = -879.3333

# Row 107 (BSM)
$ws.Range("H107").Value = 688.6842
$ws.Range("I107").Value = 694.6875
$ws.Range("J107").Value = 656.6667
$ws.Range("K107").Value = 694.6875
$ws.Range("L107").Value = 656.6667
$ws.Range("M107").Value = 1225.3125
$ws.Range("N107").Value = -4496.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 4515.8984
$ws.Range("I31").Value = 1803.2106
$ws.Range("J31").Value = 7841.129
$ws.Range("K31").Value = 1803.2106
$ws.Range("L31").Value = 7841.129
$ws.Range("M31").Value = -1508.2106
$ws.Range("N31").Value = -8431.129000000001

# Row 33 (CRP)
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

# Row 34 (CRP)
$ws.Range("H34").Value = 4515.8984
$ws.Range("I34").Value = 1803.2106
$ws.Range("J34").Value = 7841.129
$ws.Range("K34").Value = 1803.2106
$ws.Range("L34").Value = 7841.129
$ws.Range("M34").Value = -1601.2106
$ws.Range("N34").Value = -8245.129000000001

# Row 36 (CRP)
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

# Row 40 (CRP)
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# Row 125 (CRP)
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 132 (CRP)
$ws.Range("H132").Value = 737369.9
$ws.Range("I132").Value = 894264.2
$ws.Range("J132").Value = 5196.5
$ws.Range("K132").Value = 2682792.6
$ws.Range("L132").Value = 15589.5
$ws.Range("M132").Value = -2680262.6
$ws.Range("N132").Value = -20649.5

# Row 135 (CRP)
$ws.Range("H135").Value = 49679.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 49679.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 49679.5
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -59819.5

# Row 140 (CRP)
$ws.Range("H140").Value = 67299.69500000001
$ws.Range("J140").Value = 67299.69500000001
$ws.Range("L140").Value = 67299.69500000001
$ws.Range("N140").Value = -77659.69500000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 606004.2
$ws.Range("J5").Value = 1254822.2
$ws.Range("L5").Value = 3764466.6
$ws.Range("N5").Value = -3764690.6

# Row 107 (CUL)
$ws.Range("H107").Value = 469956.8
$ws.Range("J107").Value = 681161.6
$ws.Range("L107").Value = 2043484.8
$ws.Range("N107").Value = -2047324.8

# Row 108 (CUL)
$ws.Range("H108").Value = 2259.3635
$ws.Range("I108").Value = 740.375
$ws.Range("J108").Value = 6310
$ws.Range("K108").Value = 2221.125
$ws.Range("L108").Value = 18930
$ws.Range("M108").Value = 658.875
$ws.Range("N108").Value = -24690

# Row 135 (CUL)
$ws.Range("H135").Value = 606004.2
$ws.Range("J135").Value = 1254822.2
$ws.Range("L135").Value = 11293399.8
$ws.Range("N135").Value = -11298469.8

$ws = $wb.Worksheets.Item("GSM")
# Row 38 (GSM)
$ws.Range("H38").Value = 24999.5
$ws.Range("J38").Value = 24999.5
$ws.Range("L38").Value = 24999.5
$ws.Range("N38").Value = -25925.5

# Row 132 (GSM)
$ws.Range("H132").Value = 2010.1428
$ws.Range("I132").Value = 1655.44
$ws.Range("J132").Value = 4966
$ws.Range("K132").Value = 4966.32
$ws.Range("L132").Value = 14898
$ws.Range("M132").Value = -2436.32
$ws.Range("N132").Value = -19958

# Row 135 (GSM)
$ws.Range("H135").Value = 40077
$ws.Range("J135").Value = 40077
$ws.Range("L135").Value = 40077
$ws.Range("N135").Value = -50217

# Row 140 (GSM)
$ws.Range("H140").Value = 39656
$ws.Range("J140").Value = 39656
$ws.Range("L140").Value = 39656
$ws.Range("N140").Value = -50016

$ws = $wb.Worksheets.Item("LTW")
# Row 125 (LTW)
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

$ws = $wb.Worksheets.Item("WVR")
# Row 112 (WVR)
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954

# Row 137 (WVR)
$ws.Range("H137").Value = 53408.54
$ws.Range("J137").Value = 53408.54
$ws.Range("L137").Value = 53408.54
$ws.Range("N137").Value = -63608.54
